$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings that must remain TEXT (e.g. trailing
# zeros, thousand-separated "price" strings). Force text number-format before
# assigning so Excel does not silently coerce them to real numbers.

$ws.Range("D2").Value = "25.476.44"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "1.668.43"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.80"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4792"
$ws.Range("E7").Value = "  -0.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2627"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06175"
$ws.Range("E9").Value = "  +2.96%  "

$ws.Range("D10").Value = "1.667.28"
$ws.Range("E10").Value = "  +1.30%  "

$ws.Range("E11").Value = "  -2.97%  "

$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5878"
$ws.Range("E13").Value = "  -5.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.372"
$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.96"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "25.469.81"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006752"
$ws.Range("E19").Value = "  +2.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").Value = "1.881.42"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.448"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.730"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.284"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.22"
$ws.Range("E25").Value = "  +3.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.02"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.390"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.720"
$ws.Range("E28").Value = "  +3.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.90"
$ws.Range("E29").Value = "  +1.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.941"
$ws.Range("E30").Value = "  +5.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07792"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.648"
$ws.Range("E32").Value = "  +2.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9990"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("E34").Value = "  -5.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.602"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6084"
$ws.Range("E36").Value = "  +4.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9484"
$ws.Range("E37").Value = "  +1.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.591"
$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8566"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01480"
$ws.Range("E41").Value = "  -5.47%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.849"
$ws.Range("E42").Value = "  +1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.85"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3766"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.822"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1117"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.185"
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05246"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.81"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.389"
$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("E51").Value = "  +0.11%  "
